$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Week 31" column header in AF1
$ws.Range("AF1").Value = "Week 31"

# Add the new week's recorded values for the players that have one
$ws.Range("AF3").Value = 2.75
$ws.Range("AF5").Value = 7
$ws.Range("AF6").Value = 7
$ws.Range("AF7").Value = 7
$ws.Range("AF8").Value = 2

# Update the active selection to the new last-used cell
$ws.Range("AE14").Select()
